$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" value from "EF-7,EM-8" to "EF-7"
$ws.Range("B9").Value = "EF-7"
$ws.Range("C9").Value = "EF-7"

# Remove the trailing "Requisitos:" section (rows 23-24) entirely
$ws.Rows("23:24").Delete()
